$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the trailing column groups (MV_T_*, UI_*, TT_T, extra EXP_ID) that are no longer present
$ws.Columns("BV:CW").Delete()

# Update header row (BA1:BU1) with the new FORM-related column names
$ws.Range("BA1").Value = "APP"
$ws.Range("BB1").Value = "UX_A"
$ws.Range("BC1").Value = "UX_P"
$ws.Range("BD1").Value = "UX_E"
$ws.Range("BE1").Value = "UX_D"
$ws.Range("BF1").Value = "UX_S"
$ws.Range("BG1").Value = "UX_N"
$ws.Range("BH1").Value = "SUS"
$ws.Range("BI1").Value = "CL_W_MD"
$ws.Range("BJ1").Value = "CL_W_PD"
$ws.Range("BK1").Value = "CL_W_TD"
$ws.Range("BL1").Value = "CL_W_E"
$ws.Range("BM1").Value = "CL_W_P"
$ws.Range("BN1").Value = "CL_W_F"
$ws.Range("BO1").Value = "CL_MD"
$ws.Range("BP1").Value = "CL_PD"
$ws.Range("BQ1").Value = "CL_TD"
$ws.Range("BR1").Value = "CL_E"
$ws.Range("BS1").Value = "CL_P"
$ws.Range("BT1").Value = "CL_F"
$ws.Range("BU1").Value = "CL_SCORE"

# Update data rows 2-3 for the FORM columns (BA:BU)
$ws.Range("BA2").Value = 0
$ws.Range("BB2").Value = 0
$ws.Range("BC2").Value = 0
$ws.Range("BD2").Value = 0.25
$ws.Range("BE2").Value = -0.75
$ws.Range("BF2").Value = 0.25
$ws.Range("BG2").Value = 1.5
$ws.Range("BH2").Value = 90
$ws.Range("BI2").Value = 0.133
$ws.Range("BJ2").Value = 0.2
$ws.Range("BK2").Value = 0.067
$ws.Range("BL2").Value = 0.2
$ws.Range("BM2").Value = 0.067
$ws.Range("BN2").Value = 0.333
$ws.Range("BO2").Value = 0.4
$ws.Range("BP2").Value = 1.4
$ws.Range("BQ2").Value = 0.133
$ws.Range("BR2").Value = 0.8
$ws.Range("BS2").Value = 0.6
$ws.Range("BT2").Value = 0.333
$ws.Range("BU2").Value = 4.666
$ws.Range("BA3").Value = 0
$ws.Range("BB3").Value = -0.333
$ws.Range("BC3").Value = 0
$ws.Range("BD3").Value = 0.25
$ws.Range("BE3").Value = -0.25
$ws.Range("BF3").Value = -1
$ws.Range("BG3").Value = -0.25
$ws.Range("BH3").Value = 92.5
$ws.Range("BI3").Value = 0.333
$ws.Range("BJ3").Value = 0.067
$ws.Range("BK3").Value = 0.133
$ws.Range("BL3").Value = 0.2
$ws.Range("BM3").Value = 0.267
$ws.Range("BN3").Value = 0
$ws.Range("BO3").Value = 0.333
$ws.Range("BP3").Value = 0.067
$ws.Range("BQ3").Value = 0.133
$ws.Range("BR3").Value = 0.4
$ws.Range("BS3").Value = 2.133
$ws.Range("BT3").Value = 0
$ws.Range("BU3").Value = 4.066

